# "Latest updates and data fixes." - update the FoTOMRAEL input value and
# leave the workbook focused on that sheet/cell (mirrors the author's
# Excel session: FoTOMRAEL tab active, cell B3 selected).

$wb = $excel.ActiveWorkbook

# Frac of Tech Outside Modeled Region Affecting Endo Learning: 0.25 -> 0.9
$ws = $wb.Worksheets.Item("FoTOMRAEL")
$ws.Range("B2").Value = 0.9

# Switch focus to the FoTOMRAEL sheet and land the selection on B3
# (the cell just below the edited value), matching the saved view state.
$ws.Activate()
$ws.Range("B3").Select()
